# CSI 758 Project - "added the pre processing page to ppt"
#
# 1. Slide 4 ("Pre-Process Data for 3D CNN") gets its Content Placeholder
#    filled in with the pre-processing notes (and its placeholder geometry
#    gets materialized to the same box the layout already implies).
# 2. A new slide is appended ("3D CNN Using Tensor Flow") with an empty
#    content placeholder.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Slide 4 content placeholder
# ---------------------------------------------------------------------
$slide4 = $p.Slides.Item(4)
$body = $slide4.Shapes.Item(2)

# Materialize the placeholder's inherited geometry (matches slideLayout2's
# Content Placeholder box exactly: off 1154954,2603500 ext 8825659,3416300 EMU).
$body.Left = 90.94126
$body.Top = 205.0
$body.Width = 694.93384
$body.Height = 269.0

$tr = $body.TextFrame.TextRange

$para1 = "Needed to accommodate variance in number of images per Dicom"
$para2 = "Used chunking function with adjustments"
$para3 = "Approximated the number of images to get 20 groups of adjacent images"
$para4 = "Averaged the pixel values over the range of adjacent images to normalize for each group"
$dash = [char]0x2013
$para5 = "Some data fidelity loss " + $dash + " estimate that tumor size that could be observed moved from millimeters to about a centimeter in size due to averaging."
$para6 = "Result was approximately 1500 images each with 50 pixels by 50 pixels by 20 images " + $dash + " this is the basis for the common 3D model to be input to the CNN"

$full = $para1 + "`r" + $para2 + "`r" + $para3 + "`r" + $para4 + "`r" + $para5 + "`r" + $para6

# Seed with a single paragraph first so the run formatting inherits the
# placeholder's existing lang="en-US" end-paragraph mark, then grow it to
# the full text -- every paragraph created this way keeps that lang.
$tr.Text = $para1
$tr.Text = $full

# Split "Dicom" into its own run within paragraph 1 (keeps identical
# formatting -- mirrors the source deck's run boundary for that word).
$dicomIdx = $para1.IndexOf("Dicom") + 1
$dicomRun = $tr.Characters($dicomIdx, 5)
$dicomRun.Text = "Dicom"

# Bullet levels: paragraphs 3 & 4 are second-level, paragraph 5 is third-level.
$tr.Paragraphs(3, 1).IndentLevel = 2
$tr.Paragraphs(4, 1).IndentLevel = 2
$tr.Paragraphs(5, 1).IndentLevel = 3

# ---------------------------------------------------------------------
# 2. New slide: "3D CNN Using Tensor Flow"
# ---------------------------------------------------------------------
$newSlide = $p.Slides.Add($p.Slides.Count + 1, 2)
$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "3D CNN Using Tensor Flow"
